$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sicil no column becomes text ("T" + number, with row 4's number corrected 55 -> 54),
# and the role text in column B is updated for each person.

# Row 2: Sicil no "T49" with its role text
$ws.Range("A2").Value = "T49"
$ws.Range("B2").Value = "Bireysel,Yazılımcı,Yardımcı"

# Row 3: Sicil no "T50" with its role text
$ws.Range("A3").Value = "T50"
$ws.Range("B3").Value = "Araştırmacı,Yazılımcı"

# Row 4: Sicil no "T54" with its role text
$ws.Range("A4").Value = "T54"
$ws.Range("B4").Value = "Yardımcı"
